$d = $word.ActiveDocument

# --- Locate the anchor paragraphs --------------------------------------
# The "<별첨>" paragraph is the only paragraph whose text starts with "<".
# The paragraph immediately before it is the page-break paragraph that
# needs to be replaced with three new bulleted list items.
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.Length -gt 0 -and $t.Substring(0,1) -eq "<") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "could not locate the <별첨> paragraph"
}

$pageBreakPara = $d.Paragraphs($anchorIndex - 1)
$attachPara    = $d.Paragraphs($anchorIndex)
$picturePara   = $d.Paragraphs($anchorIndex + 1)

$wNs     = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$bulletPPr = '<w:pPr><w:pStyle w:val="a6"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr></w:pPr>'

# --- Replace the page-break paragraph with three new list paragraphs ---
$bullet1 = '<w:p ' + $wNs + '>' + $bulletPPr + `
    '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t>void printColorStrip(int colorNum)</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> 을 통해 0, 1, 2, 3(흰색 스트립, 녹색, 청색, 보라색)의 컬러 스티커를 출력 할 수 있도록 하였습니다. </w:t></w:r>' + `
    '</w:p>'

$bullet2 = '<w:p ' + $wNs + '>' + $bulletPPr + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t>9999 시간에 대해 All Day Long 텍스트의 출력 및 최우선 정렬이 가능합니다.</w:t></w:r>' + `
    '</w:p>'

$bullet3 = '<w:p ' + $wNs + '>' + $bulletPPr + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>(버그 수정) 공휴일로 고정된 일정이 계속해서 증식하는 문제를 해결했습니다.</w:t></w:r>' + `
    '</w:p>'

$pageBreakPara.Range.InsertXML($bullet1 + $bullet2 + $bullet3)

# Paragraph indices shift by +2 after the insert above (1 removed, 3 added).
$attachIndex = $anchorIndex + 2
$attachPara  = $d.Paragraphs($attachIndex)
$picturePara = $d.Paragraphs($attachIndex + 1)

# --- "<별첨>" paragraph keeps its text/run formatting but no longer ----
# --- carries the (now relocated) lastRenderedPageBreak marker. ---------
$attachXml = '<w:p ' + $wNs + '><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>&lt;별첨&gt;</w:t></w:r></w:p>'
$attachPara.Range.InsertXML($attachXml)

# --- The picture paragraph is untouched content-wise; only the ---------
# --- drawing's wp14:editId attribute is refreshed (as Word does on ----
# --- every re-save that touches the surrounding flow). -----------------
$pictureXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:drawing>' + `
    '<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="35E3820C" wp14:editId="74AB0302">' + `
    '<wp:extent cx="6629400" cy="3803650"/>' + `
    '<wp:effectExtent l="0" t="0" r="0" b="6350"/>' + `
    '<wp:docPr id="948048947" name="그림 1"/>' + `
    '<wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
    '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
    '<pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr>' + `
    '<pic:blipFill><a:blip r:embed="rId7" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
    '<pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="6629400" cy="3803650"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr>' + `
    '</pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'
$picturePara.Range.InsertXML($pictureXml)

Write-Output "done"
